$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded. Insert a row at position 62
# (pushing the existing row 62 and everything below it down by one row,
# growing the table from A1:R147 to A1:R148) and populate it with the new
# record's data.
$ws.Rows("62:62").Insert()

$ws.Range("A62").Value2 = 3
$ws.Range("B62").Value2 = "Femacal de La Calera"
$ws.Range("C62").Value2 = "Coquimbo"
$ws.Range("D62").Value2 = 44740
$ws.Range("E62").Value2 = 5
$ws.Range("F62").Value2 = 100112026
$ws.Range("G62").Value2 = "Haba"
$ws.Range("H62").Value2 = "Sin especificar"
$ws.Range("I62").Value2 = "Primera"
$ws.Range("J62").Value2 = 103
$ws.Range("K62").Value2 = 21000
$ws.Range("L62").Value2 = 22000
$ws.Range("M62").Value2 = 21534
$ws.Range("N62").Value2 = "`$/saco 25 kilos"
$ws.Range("O62").Value2 = "Provincia de Limarí"
$ws.Range("P62").Value2 = 861
$ws.Range("Q62").Value2 = 25
$ws.Range("R62").Value2 = "Hortaliza"
